$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Update the ACOMPANHAMENTO column (O2:O10) with new CIP values
$ws.Range("O2").Value = "CIP352"
$ws.Range("O3").Value = "CIP353"
$ws.Range("O4").Value = "CIP354"
$ws.Range("O5").Value = "CIP355"
$ws.Range("O6").Value = "CIP356"
$ws.Range("O7").Value = "CIP357"
$ws.Range("O8").Value = "CIP358"
$ws.Range("O9").Value = "CIP359"
$ws.Range("O10").Value = "CIP360"

# Update the selected cell in the sheet view
$ws.Range("M13").Select()
